$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 182, shifting the existing
# rows 182-189 down to 184-191 (the data itself does not change, only its
# position moves down by two rows).
$ws.Rows.Item(182).Insert()
$ws.Rows.Item(182).Insert()

# Fill in the new row 182 ("Primera" quality) with the latest week's data.
$ws.Cells.Item(182, 1).Value = 8
$ws.Cells.Item(182, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44516
$ws.Cells.Item(182, 4).NumberFormat = $ws.Cells.Item(184, 4).NumberFormat
$ws.Cells.Item(182, 5).Value = 4
$ws.Cells.Item(182, 6).Value = 100114014
$ws.Cells.Item(182, 7).Value = "Betarraga"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 3060
$ws.Cells.Item(182, 11).Value = 450
$ws.Cells.Item(182, 12).Value = 500
$ws.Cells.Item(182, 13).Value = 475
$ws.Cells.Item(182, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(182, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(182, 16).Value = 158
$ws.Cells.Item(182, 17).Value = 3
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# Fill in the new row 183 ("Segunda" quality) with the latest week's data.
$ws.Cells.Item(183, 1).Value = 8
$ws.Cells.Item(183, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44516
$ws.Cells.Item(183, 4).NumberFormat = $ws.Cells.Item(185, 4).NumberFormat
$ws.Cells.Item(183, 5).Value = 4
$ws.Cells.Item(183, 6).Value = 100114014
$ws.Cells.Item(183, 7).Value = "Betarraga"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Segunda"
$ws.Cells.Item(183, 10).Value = 1480
$ws.Cells.Item(183, 11).Value = 350
$ws.Cells.Item(183, 12).Value = 400
$ws.Cells.Item(183, 13).Value = 375
$ws.Cells.Item(183, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(183, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(183, 16).Value = 125
$ws.Cells.Item(183, 17).Value = 3
$ws.Cells.Item(183, 18).Value = "Hortaliza"
